# Generate Report for Archive
#
# 1. Update status text "Ready for handoff" -> "In Translation" on every
#    sheet/cell where it appears (Overview!E2:F2, zh-cn!C2, de-de!C2).
# 2. Narrow the "Status" column(s) from ~17.22 chars to ~13.41 chars on
#    every sheet (Overview columns E & F, zh-cn column C, de-de column C).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# The target OOXML column width (13.4101845877511 chars) corresponds to a
# COM ColumnWidth of about 12.58 characters (width = ColumnWidth + 5/6).
$newColumnWidth = 12.576851254417766

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
if ($wsOverview.Range("E2").Value2 -eq $oldStatus) {
    $wsOverview.Range("E2").Value = $newStatus
}
if ($wsOverview.Range("F2").Value2 -eq $oldStatus) {
    $wsOverview.Range("F2").Value = $newStatus
}
$wsOverview.Range("E1:F1").ColumnWidth = $newColumnWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
if ($wsZhCn.Range("C2").Value2 -eq $oldStatus) {
    $wsZhCn.Range("C2").Value = $newStatus
}
$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
if ($wsDeDe.Range("C2").Value2 -eq $oldStatus) {
    $wsDeDe.Range("C2").Value = $newStatus
}
$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
